$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 14.55576705260431
$ws.Range("C2").Value2 = 9.715555350721269
$ws.Range("D2").Value2 = 5.967597785388032
$ws.Range("E2").Value2 = 16.53464518793436
$ws.Range("G2").Value2 = 3.605956699978702
$ws.Range("N2").Value2 = 16.15533835233212
$ws.Range("O2").Value2 = 20.97852255951449
$ws.Range("B3").Value2 = 13.8873474131362
$ws.Range("C3").Value2 = 9.131014552356973
$ws.Range("D3").Value2 = 5.846138528907759
$ws.Range("E3").Value2 = 15.59184432214095
$ws.Range("G3").Value2 = 3.608895072387722
$ws.Range("N3").Value2 = 16.2213662345871
$ws.Range("O3").Value2 = 20.91500924417518
$ws.Range("B4").Value2 = 13.46258790654218
$ws.Range("C4").Value2 = 8.750897233833163
$ws.Range("D4").Value2 = 5.772161414271669
$ws.Range("E4").Value2 = 14.9882413140507
$ws.Range("G4").Value2 = 3.610791680355224
$ws.Range("N4").Value2 = 16.26382801139703
$ws.Range("O4").Value2 = 20.88358147671162
$ws.Range("B5").Value2 = 13.28613472166723
$ws.Range("C5").Value2 = 8.59069028053313
$ws.Range("D5").Value2 = 5.742212253355198
$ws.Range("E5").Value2 = 14.73632155765062
$ws.Range("G5").Value2 = 3.611587895181681
$ws.Range("N5").Value2 = 16.28161604525777
$ws.Range("O5").Value2 = 20.87268113028909
$ws.Range("B6").Value2 = 13.25663978220878
$ws.Range("C6").Value2 = 8.563768524901469
$ws.Range("D6").Value2 = 5.737252486834846
$ws.Range("E6").Value2 = 14.69414023316585
$ws.Range("G6").Value2 = 3.611721517808098
$ws.Range("N6").Value2 = 16.28459904596025
$ws.Range("O6").Value2 = 20.87098633863554
$ws.Range("B7").Value2 = 13.46022144977527
$ws.Range("C7").Value2 = 8.748758069028357
$ws.Range("D7").Value2 = 5.771756649678003
$ws.Range("E7").Value2 = 14.98486751869946
$ws.Range("G7").Value2 = 3.610802323802448
$ws.Range("N7").Value2 = 16.26406594315362
$ws.Range("O7").Value2 = 20.88342674846652
$ws.Range("B8").Value2 = 14.32841151812921
$ws.Range("C8").Value2 = 9.518420395404013
$ws.Range("D8").Value2 = 5.925621914877119
$ws.Range("E8").Value2 = 16.21484378950503
$ws.Range("G8").Value2 = 3.606950718937489
$ws.Range("N8").Value2 = 16.17770716320013
$ws.Range("O8").Value2 = 20.95505488076402
$ws.Range("B9").Value2 = 15.90790783178707
$ws.Range("C9").Value2 = 10.85871366111224
$ws.Range("D9").Value2 = 6.230102655258841
$ws.Range("E9").Value2 = 18.5217318491883
$ws.Range("G9").Value2 = 3.600127145057569
$ws.Range("N9").Value2 = 16.02352121577687
$ws.Range("O9").Value2 = 21.15530716330676
$ws.Range("B10").Value2 = 16.9831606796559
$ws.Range("C10").Value2 = 11.73986517357016
$ws.Range("D10").Value2 = 6.452879297605894
$ws.Range("E10").Value2 = 20.17637832381546
$ws.Range("G10").Value2 = 3.595552896254112
$ws.Range("N10").Value2 = 15.91938070916921
$ws.Range("O10").Value2 = 21.33833388838116
$ws.Range("B11").Value2 = 17.45218223088731
$ws.Range("C11").Value2 = 12.11822433975327
$ws.Range("D11").Value2 = 6.553506457121886
$ws.Range("E11").Value2 = 20.8870772796232
$ws.Range("G11").Value2 = 3.593566077809811
$ws.Range("N11").Value2 = 15.87396714250416
$ws.Range("O11").Value2 = 21.4292197672453
$ws.Range("B12").Value2 = 17.62678780445525
$ws.Range("C12").Value2 = 12.25826899067452
$ws.Range("D12").Value2 = 6.591468065513268
$ws.Range("E12").Value2 = 21.15019811576933
$ws.Range("G12").Value2 = 3.592827150239774
$ws.Range("N12").Value2 = 15.85705049115966
$ws.Range("O12").Value2 = 21.46471386793368
$ws.Range("B13").Value2 = 17.58931847482856
$ws.Range("C13").Value2 = 12.22825155792136
$ws.Range("D13").Value2 = 6.583299371927011
$ws.Range("E13").Value2 = 21.09379686970399
$ws.Range("G13").Value2 = 3.592985695212442
$ws.Range("N13").Value2 = 15.86068134281278
$ws.Range("O13").Value2 = 21.45702198096397
$ws.Range("B14").Value2 = 17.46660788050926
$ws.Range("C14").Value2 = 12.12981065126924
$ws.Range("D14").Value2 = 6.556632687666879
$ws.Range("E14").Value2 = 20.90884454568253
$ws.Range("G14").Value2 = 3.593505016967777
$ws.Range("N14").Value2 = 15.87256978629588
$ws.Range("O14").Value2 = 21.43211840941423
$ws.Range("B15").Value2 = 17.39105010527932
$ws.Range("C15").Value2 = 12.06909196007603
$ws.Range("D15").Value2 = 6.54027868389707
$ws.Range("E15").Value2 = 20.79477486543317
$ws.Range("G15").Value2 = 3.593824864002857
$ws.Range("N15").Value2 = 15.87988828592248
$ws.Range("O15").Value2 = 21.4170039803921
$ws.Range("B16").Value2 = 16.95209439667316
$ws.Range("C16").Value2 = 11.71468615670428
$ws.Range("D16").Value2 = 6.446285139635616
$ws.Range("E16").Value2 = 20.12909018617719
$ws.Range("G16").Value2 = 3.595684624455787
$ws.Range("N16").Value2 = 15.92238791927816
$ws.Range("O16").Value2 = 21.33254616824333
$ws.Range("B17").Value2 = 16.67757273254012
$ws.Range("C17").Value2 = 11.49151153214598
$ws.Range("D17").Value2 = 6.388410186511368
$ws.Range("E17").Value2 = 19.7099800922242
$ws.Range("G17").Value2 = 3.596849550810268
$ws.Range("N17").Value2 = 15.94896116516624
$ws.Range("O17").Value2 = 21.28267450753444
$ws.Range("B18").Value2 = 16.51778698047372
$ws.Range("C18").Value2 = 11.36103110470045
$ws.Range("D18").Value2 = 6.355057045568628
$ws.Range("E18").Value2 = 19.46496296017253
$ws.Range("G18").Value2 = 3.597528441218091
$ws.Range("N18").Value2 = 15.9644300065965
$ws.Range("O18").Value2 = 21.25470875299639
$ws.Range("B19").Value2 = 16.4633656378119
$ws.Range("C19").Value2 = 11.31648921593615
$ws.Range("D19").Value2 = 6.343754364932163
$ws.Range("E19").Value2 = 19.38132382938809
$ws.Range("G19").Value2 = 3.597759825362075
$ws.Range("N19").Value2 = 15.96969924199637
$ws.Range("O19").Value2 = 21.24536407551974
$ws.Range("B20").Value2 = 16.70699230446739
$ws.Range("C20").Value2 = 11.5154878808408
$ws.Range("D20").Value2 = 6.394578104665507
$ws.Range("E20").Value2 = 19.75500426050749
$ws.Range("G20").Value2 = 3.596724626505552
$ws.Range("N20").Value2 = 15.94611330156321
$ws.Range("O20").Value2 = 21.28790913661439
$ws.Range("B21").Value2 = 17.5027332354797
$ws.Range("C21").Value2 = 12.15881281925251
$ws.Range("D21").Value2 = 6.564469553930527
$ws.Range("E21").Value2 = 20.96333223222608
$ws.Range("G21").Value2 = 3.593352115558256
$ws.Range("N21").Value2 = 15.86907026140341
$ws.Range("O21").Value2 = 21.43940411027432
$ws.Range("B22").Value2 = 18.00525519707561
$ws.Range("C22").Value2 = 12.56042799264554
$ws.Range("D22").Value2 = 6.674648465542629
$ws.Range("E22").Value2 = 21.71805637447925
$ws.Range("G22").Value2 = 3.591226269702703
$ws.Range("N22").Value2 = 15.82035238828039
$ws.Range("O22").Value2 = 21.54468534580228
$ws.Range("B23").Value2 = 17.73868648240737
$ws.Range("C23").Value2 = 12.34780079954224
$ws.Range("D23").Value2 = 6.615935039677584
$ws.Range("E23").Value2 = 21.31843570215009
$ws.Range("G23").Value2 = 3.592353738231882
$ws.Range("N23").Value2 = 15.84620497154439
$ws.Range("O23").Value2 = 21.48792792519243
$ws.Range("B24").Value2 = 16.69369781113884
$ws.Range("C24").Value2 = 11.50465494639514
$ws.Range("D24").Value2 = 6.391789835125721
$ws.Range("E24").Value2 = 19.73466149922766
$ws.Range("G24").Value2 = 3.596781076248702
$ws.Range("N24").Value2 = 15.94740022396785
$ws.Range("O24").Value2 = 21.28554035982094
$ws.Range("B25").Value2 = 15.49494189167222
$ws.Range("C25").Value2 = 10.5142977182947
$ws.Range("D25").Value2 = 6.147712538416225
$ws.Range("E25").Value2 = 17.87490713031399
$ws.Range("G25").Value2 = 3.601895598056553
$ws.Range("N25").Value2 = 16.0636200984109
$ws.Range("O25").Value2 = 21.09477845081138

Write-Host "Applied 380 kV case updates: loading_percent values refreshed"